# Insert a new data row just above row 446 (shifts rows 446:529 down to 447:530)
# and populate it with a new record. This mirrors: dimension A1:R529 -> A1:R530.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(446).Insert()

$ws.Cells.Item(446, 1).Value = 4
$ws.Cells.Item(446, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(446, 3).Value = "Los Lagos"
$ws.Cells.Item(446, 4).Value = 45209
$ws.Cells.Item(446, 5).Value = 10
$ws.Cells.Item(446, 6).Value = 100112003
$ws.Cells.Item(446, 7).Value = "Ajo"
$ws.Cells.Item(446, 8).Value = "Chino"
$ws.Cells.Item(446, 9).Value = "Primera"
$ws.Cells.Item(446, 10).Value = 240
$ws.Cells.Item(446, 11).Value = 24000
$ws.Cells.Item(446, 12).Value = 26000
$ws.Cells.Item(446, 13).Value = 25000
$ws.Cells.Item(446, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(446, 15).Value = "China"
$ws.Cells.Item(446, 16).Value = 2500
$ws.Cells.Item(446, 17).Value = 10
$ws.Cells.Item(446, 18).Value = "Hortaliza"
